$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("EntryPositionsManager", 2),
    @("ExitPositionsManager", 2),
    @("IncentivesVault", 6),
    @("InterestRatesManager", 6),
    @("IAaveDistributionManager", 0),
    @("IAaveIncentivesController", 0),
    @("IAToken", 2),
    @("IERC20", 0),
    @("ILendingPool", 2),
    @("ILendingPoolAddressesProvider", 0),
    @("IPriceOracleGetter", 0),
    @("IScaledBalanceToken", 0),
    @("IVariableDebtToken", 1),
    @("IEntryPositionsManager", 0),
    @("IExitPositionsManager", 0),
    @("IGetterUnderlyingAsset", 0),
    @("IIncentivesVault", 1),
    @("IInterestRatesManager", 0),
    @("IMorpho", 8),
    @("IOracle", 0),
    @("IRewardsManager", 1),
    @("ILido", 0),
    @("IndexesLens", 3),
    @("ILens", 3),
    @("Lens", 1),
    @("LensStorage", 11),
    @("MarketsLens", 1),
    @("RatesLens", 2),
    @("UsersLens", 2),
    @("DataTypes", 0),
    @("Errors", 0),
    @("ReserveConfiguration", 2),
    @("UserConfiguration", 1),
    @("InterestRatesModel", 4),
    @("Types", 0),
    @("MatchingEngine", 1),
    @("Morpho", 1),
    @("MorphoGovernance", 2),
    @("MorphoStorage", 10),
    @("MorphoUtils", 9),
    @("PositionsManagerUtils", 3),
    @("RewardsDistributor", 3),
    @("FakeToken", 1),
    @("IncentivesVault", 6),
    @("InterestRatesManager", 4),
    @("ICEth", 0),
    @("IComptroller", 0),
    @("IInterestRateModel", 0),
    @("ICToken", 0),
    @("ICEther", 0),
    @("ICompoundOracle", 0),
    @("IIncentivesVault", 1),
    @("IInterestRatesManager", 0),
    @("IMorpho", 5),
    @("IOracle", 0),
    @("IPositionsManager", 0),
    @("IRewardsManager", 1),
    @("IWETH", 0),
    @("IndexesLens", 2),
    @("ILens", 3),
    @("Lens", 1),
    @("LensStorage", 9),
    @("MarketsLens", 1),
    @("RatesLens", 1),
    @("RewardsLens", 1),
    @("UsersLens", 1),
    @("CompoundMath", 0),
    @("InterestRatesModel", 4),
    @("Types", 0),
    @("MatchingEngine", 1),
    @("Morpho", 1),
    @("MorphoGovernance", 1),
    @("MorphoStorage", 9),
    @("MorphoUtils", 5),
    @("PositionsManager", 3),
    @("RewardsManager", 4)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = 0
    $r++
}
